$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing default tier values that changed
$ws.Range("C17").Value = 2
$ws.Range("C19").Value = 2

# Add two new rows for missing types: Chiropractic and Dialysis (both "health clinic", tier 2)
$ws.Range("A20").Value = "Chiropractic"
$ws.Range("B20").Value = "health clinic"
$ws.Range("C20").Value = 2

$ws.Range("A21").Value = "Dialysis"
$ws.Range("B21").Value = "health clinic"
$ws.Range("C21").Value = 2

# The stale "D1:D1048576" selection left over from a previous edit is no
# longer relevant now that the table has grown - reset the active
# selection back to the top-left cell.
$ws.Range("A1").Select()
